$d = $word.ActiveDocument

$pairs = @(
    @("2024-03-19 Tuesday", "2024-03-20 Wednesday"),
    @("45+47=", "2+71="),
    @("1+11=", "55-53="),
    @("54-20=", "21+18="),
    @("98-51=", "21-21="),
    @("87-19=", "13+63="),
    @("69-12=", "40+38="),
    @("91-7=", "47+36="),
    @("95-91=", "71-3="),
    @("9+13=", "65-18="),
    @("96-24=", "1+58="),
    @("43-24=", "50-7="),
    @("1+61=", "62-59="),
    @("84-52=", "55-29="),
    @("66-31=", "99-68="),
    @("46+27=", "41+47="),
    @("84-72=", "95-81="),
    @("36+41=", "64+29="),
    @("48-28=", "86-75="),
    @("4+34=", "66-24="),
    @("10+64=", "49-11="),
    @("34+52=", "46-5="),
    @("14+12=", "44+21="),
    @("2+33=", "78+10="),
    @("70-21=", "61-39="),
    @("34+57=", "22+57="),
    @("97-70=", "53-1="),
    @("43-31=", "61-59="),
    @("69-45=", "7+7="),
    @("85+3=", "52-8="),
    @("54+30=", "31+5="),
    @("48+20=", "93-46="),
    @("74-46=", "59-4="),
    @("1+42=", "3+65="),
    @("88-15=", "19+65="),
    @("66-37=", "72+6="),
    @("14+31=", "1+60="),
    @("42-8=", "81-9="),
    @("42-24=", "95-78="),
    @("39+22=", "1+67="),
    @("5+73=", "3+0="),
    @("55+16=", "63-6="),
    @("22+59=", "67-15="),
    @("71-47=", "66+9="),
    @("28+18=", "81-62="),
    @("89-2=", "67+4="),
    @("3+72=", "77-75="),
    @("63+23=", "99-44="),
    @("70+3=", "45+14="),
    @("84-17=", "9+3="),
    @("19+7=", "61-28="),
    @("6+20=", "86-62="),
    @("98-91=", "80-56="),
    @("7+55=", "83-5="),
    @("21+10=", "60+17="),
    @("28+17=", "9+25="),
    @("7+15=", "84-10="),
    @("17+36=", "64-46="),
    @("70+27=", "82+16="),
    @("69-49=", "0+43="),
    @("70-5=", "25+5="),
    @("44-36=", "57-10="),
    @("75-54=", "0+29="),
    @("69-27=", "93-92="),
    @("47-4=", "40+44="),
    @("39+25=", "60-25="),
    @("51-45=", "78+3="),
    @("4+10=", "38-0="),
    @("21+72=", "40+51="),
    @("56+17=", "71-66="),
    @("25+22=", "75+3="),
    @("87-23=", "58-0="),
    @("72-11=", "80-70="),
    @("44+18=", "29+50="),
    @("0+39=", "63-57="),
    @("41-37=", "46+51="),
    @("35-1=", "6-3="),
    @("15+2=", "84-4="),
    @("13+70=", "59+28="),
    @("61+36=", "21-15="),
    @("17+9=", "66-35="),
    @("48-38=", "55-34="),
    @("56-49=", "81-54="),
    @("49+21=", "60-30="),
    @("84-68=", "7+67="),
    @("57-22=", "51+30="),
    @("90-73=", "63-4="),
    @("0+45=", "82-24="),
    @("66-65=", "67-16="),
    @("18+37=", "47-22="),
    @("67-18=", "6+5="),
    @("54+22=", "1+0="),
    @("24+74=", "15+45="),
    @("34-0=", "93-9="),
    @("7+75=", "80-22="),
    @("52+47=", "81-72="),
    @("43-3=", "48-9="),
    @("5+71=", "75-70="),
    @("29+34=", "32+5="),
    @("32-27=", "67-40="),
    @("91-29=", "71-32="),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

Write-Output "Done"
